$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Variables")
$ws2 = $wb.Worksheets.Item("MeshSizes")
$ws3 = $wb.Worksheets.Item("Volume")

# ---------------------------------------------------------------------------
# Sheet "Variables": rebuild as a describe()-style table with an extra
# "Unnamed: 0" index column and a 5th data column ("fine").
# ---------------------------------------------------------------------------
$ws1.Range("A1").Value = "Unnamed: 0"
$ws1.Range("B1").Value = "variable"
$ws1.Range("C1").Value = "coarse"
$ws1.Range("D1").Value = "medium"
$ws1.Range("E1").Value = "fine"

# E1 is a brand new header cell; give it the same bold/centered/bordered
# look as the rest of the header row by copying the formatting from D1.
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "Mean"
$ws1.Range("C2").Value = 0.6306836866502078
$ws1.Range("D2").Value = 0.6187786730549253
$ws1.Range("E2").Value = 0.6126318529329569

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "Standard Deviation"
$ws1.Range("C3").Value = 0.208289031674727
$ws1.Range("D3").Value = 0.2107123600455134
$ws1.Range("E3").Value = 0.2113512510940114

$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "Variance"
$ws1.Range("C4").Value = 0.04338432071599541
$ws1.Range("D4").Value = 0.04439969867595008
$ws1.Range("E4").Value = 0.04466935133900386

$ws1.Range("A5").Value = 3
$ws1.Range("B5").Value = "Coefficient of Variation"
$ws1.Range("C5").Value = 33.02591078279293
$ws1.Range("D5").Value = 34.05294481227374
$ws1.Range("E5").Value = 34.49890012773145

# ---------------------------------------------------------------------------
# Sheet "MeshSizes": headers keep their original text; only the data row
# changes.
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = 8892
$ws2.Range("B2").Value = 62711
$ws2.Range("C2").Value = 470421

# ---------------------------------------------------------------------------
# Sheet "Volume": header keeps its original text; data value updated.
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = 0.44
